# Apply the crypto price/volume refresh captured in the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '43.803.82'
$ws.Cells.Item(2, 5).Value = '  +4.83%  '
$ws.Cells.Item(3, 4).Value = '2.282.38'
$ws.Cells.Item(3, 5).Value = '  +2.40%  '
$ws.Cells.Item(4, 5).Value = '  +0.10%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '231.50'
$ws.Cells.Item(5, 5).Value = '  -0.13%  '
$ws.Cells.Item(6, 4).Value = '0.622'
$ws.Cells.Item(6, 5).Value = '  -0.27%  '
$ws.Cells.Item(7, 4).Value = '64.74'
$ws.Cells.Item(7, 5).Value = '  +6.82%  '
$ws.Cells.Item(8, 5).Value = '  +0.07%  '
$ws.Cells.Item(9, 5).Value = '  +4.72%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0970'
$ws.Cells.Item(10, 5).Value = '  +7.65%  '
$ws.Cells.Item(11, 4).Value = '57.98'
$ws.Cells.Item(11, 5).Value = '  -0.56%  '
$ws.Cells.Item(12, 4).Value = '26.32'
$ws.Cells.Item(12, 5).Value = '  +15.70%  '
$ws.Cells.Item(13, 5).Value = '  +0.37%  '
$ws.Cells.Item(14, 4).Value = '2.622.51'
$ws.Cells.Item(14, 5).Value = '  +2.45%  '
$ws.Cells.Item(15, 4).Value = '15.79'
$ws.Cells.Item(15, 5).Value = '  +1.08%  '
$ws.Cells.Item(16, 5).Value = '  +5.53%  '
$ws.Cells.Item(17, 4).Value = '0.819'
$ws.Cells.Item(17, 5).Value = '  +2.26%  '
$ws.Cells.Item(18, 4).Value = '2.296.39'
$ws.Cells.Item(18, 5).Value = '  +2.50%  '
$ws.Cells.Item(19, 4).Value = '43.689.86'
$ws.Cells.Item(19, 5).Value = '  +4.69%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0951'
$ws.Cells.Item(20, 5).Value = '  +4.79%  '
$ws.Cells.Item(21, 4).Value = '73.43'
$ws.Cells.Item(21, 5).Value = '  +1.25%  '
$ws.Cells.Item(22, 4).Value = '6.18'
$ws.Cells.Item(22, 5).Value = '  +0.89%  '
$ws.Cells.Item(23, 4).Value = '250.37'
$ws.Cells.Item(23, 5).Value = '  +0.98%  '
$ws.Cells.Item(24, 5).Value = '  +0.11%  '
$ws.Cells.Item(25, 4).Value = '2.54'
$ws.Cells.Item(25, 5).Value = '  +6.16%  '
$ws.Cells.Item(26, 5).Value = '  +1.11%  '
$ws.Cells.Item(27, 4).Value = '9.99'
$ws.Cells.Item(27, 5).Value = '  +3.37%  '
$ws.Cells.Item(28, 4).Value = '171.89'
$ws.Cells.Item(28, 5).Value = '  +1.59%  '
$ws.Cells.Item(29, 5).Value = '  -2.22%  '
$ws.Cells.Item(30, 4).Value = '20.59'
$ws.Cells.Item(30, 5).Value = '  +3.35%  '
$ws.Cells.Item(31, 4).Value = '1.44'
$ws.Cells.Item(31, 5).Value = '  +3.00%  '
$ws.Cells.Item(32, 5).Value = '  +5.06%  '
$ws.Cells.Item(33, 4).Value = '0.122'
$ws.Cells.Item(33, 5).Value = '  +0.52%  '
$ws.Cells.Item(34, 4).Value = '5.27'
$ws.Cells.Item(34, 5).Value = '  +5.29%  '
$ws.Cells.Item(35, 5).Value = '  +6.66%  '
$ws.Cells.Item(36, 4).Value = '4.74'
$ws.Cells.Item(36, 5).Value = '  +0.65%  '
$ws.Cells.Item(37, 4).Value = '6.83'
$ws.Cells.Item(37, 5).Value = '  +4.05%  '
$ws.Cells.Item(38, 4).Value = '3.78'
$ws.Cells.Item(38, 5).Value = '  +4.46%  '
$ws.Cells.Item(39, 4).Value = '2.35'
$ws.Cells.Item(39, 5).Value = '  -1.43%  '
$ws.Cells.Item(40, 5).Value = '  +3.55%  '
$ws.Cells.Item(41, 5).Value = '  +0.16%  '
$ws.Cells.Item(42, 4).Value = '11.05'
$ws.Cells.Item(42, 5).Value = '  +28.16%  '
$ws.Cells.Item(43, 2).Value = 'FTXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(43, 4).Value = '4.74'
$ws.Cells.Item(43, 5).Value = '  +6.04%  '
$ws.Cells.Item(44, 2).Value = 'TerraClassic'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Cells.Item(44, 4).Value = '0.000228'
$ws.Cells.Item(44, 5).Value = '  -2.92%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '8.50'
$ws.Cells.Item(45, 5).Value = '  -0.82%  '
$ws.Cells.Item(46, 5).Value = '  +0.37%  '
$ws.Cells.Item(47, 4).Value = '0.0965'
$ws.Cells.Item(47, 5).Value = '  +0.63%  '
$ws.Cells.Item(48, 4).Value = '98.14'
$ws.Cells.Item(48, 5).Value = '  -0.51%  '
$ws.Cells.Item(49, 4).Value = '1.488.86'
$ws.Cells.Item(49, 5).Value = '  +1.21%  '
$ws.Cells.Item(50, 4).Value = '16.92'
$ws.Cells.Item(50, 5).Value = '  +2.03%  '
$ws.Cells.Item(51, 4).Value = '2.34'
$ws.Cells.Item(51, 5).Value = '  +1.64%  '
